# Parameter.xlsx - "Hinzufügen der Quellen" edit
# Adds source/comment references (column E "Quelle", some column D "Kommentar")
# across several rows, marks a few "ask Tim" placeholders in red, and turns
# the roof-pitch source into a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8 - Giebeldachfläche: add a clickable source hyperlink in E8
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add(
    $ws.Range("E8"),
    "https://www.kern-haus.de/ratgeber/baulexikon/satteldach/",
    "",
    "",
    "https://www.kern-haus.de/ratgeber/baulexikon/satteldach/"
) | Out-Null

# ---------------------------------------------------------------------
# Rows 20-22 - Luftdichte/Wärmekapazität: sources still unclear -> "Tim fragen"
# (marked in red as an open TODO)
# ---------------------------------------------------------------------
foreach ($r in 20, 21, 22) {
    $cell = $ws.Range("E$r")
    $cell.Value = "Tim fragen"
    $cell.Font.Color = 255
}

# ---------------------------------------------------------------------
# Rows 32-33 - tolerabler Druckverlust: max-value comment + lecture source
# ---------------------------------------------------------------------
foreach ($r in 32, 33) {
    $ws.Range("D$r").Value = "Maximalwert"
    $ws.Range("E$r").Value = "Vorlesung 04 Folie 9"
}

# ---------------------------------------------------------------------
# Rows 34-37 - Vor-/Rücklauftemperaturen Winter/Sommer: lecture source
# ---------------------------------------------------------------------
foreach ($r in 34, 35, 36, 37) {
    $ws.Range("E$r").Value = "Vorlesung 03 Folie 24 Beründung siehe Doku_all"
}

# ---------------------------------------------------------------------
# Row 38 - Bodentemperatur Winter: still open -> "Tim fragen" (red)
# ---------------------------------------------------------------------
$cell = $ws.Range("E38")
$cell.Value = "Tim fragen"
$cell.Font.Color = 255

# ---------------------------------------------------------------------
# Row 39 - Bodentemperatur Sommer: comment + still open -> "Tim fragen" (red)
# ---------------------------------------------------------------------
$ws.Range("D39").Value = "17,5°C Bodentemperatur"
$cell = $ws.Range("E39")
$cell.Value = "Tim fragen"
$cell.Font.Color = 255

# ---------------------------------------------------------------------
# Row 41 - Strömungsgeschwindigkeit: comment + lecture source
# ---------------------------------------------------------------------
$ws.Range("D41").Value = "Maximale Stromungsgeschwindigkeit auf Grund von Geräuschemissionen"
$ws.Range("E41").Value = "Vorlesung 02 Folie 10"

# ---------------------------------------------------------------------
# Rows 42-45 - Dichte/Viskosität Vor-/Rücklauf: peacesoftware.de source
# (44/45 also get the same "4 bar" comment already used on 42/43)
# ---------------------------------------------------------------------
$ws.Range("E42").Value = "http://www.peacesoftware.de/einigewerte/wasser_dampf.html"
$ws.Range("E43").Value = "http://www.peacesoftware.de/einigewerte/wasser_dampf.html"

$ws.Range("D44").Value = "4 bar, 75°C peacesoftware.de Auf Druck und Temperatur anpassen!"
$ws.Range("E44").Value = "http://www.peacesoftware.de/einigewerte/wasser_dampf.html"

$ws.Range("D45").Value = "4 bar, 55°C peacesoftware.de Auf Druck und Temperatur anpassen!"
$ws.Range("E45").Value = "http://www.peacesoftware.de/einigewerte/wasser_dampf.html"

# ---------------------------------------------------------------------
# Row 47 - Initiale Dämmung: comment
# ---------------------------------------------------------------------
$ws.Range("D47").Value = "Wir beginnen mit der günstigsten Insolationsstufe"

# ---------------------------------------------------------------------
# Row 48 - Rohrrauheit k: comment + lecture source
# ---------------------------------------------------------------------
$ws.Range("D48").Value = "Annahme Rohrrauheit"
$ws.Range("E48").Value = "Vorlesung 00 Hinweise zur Bearbeitung Folie 10 Punkt 7-9 Annahme Rohrrauheit"

# ---------------------------------------------------------------------
# Row 49 - Untergrenze Hydraulisch glatt: comment + source
# ---------------------------------------------------------------------
$ws.Range("D49").Value = "Reynolds-Zahl"
$ws.Range("E49").Value = "Planungshandbuch Seite 128"

# ---------------------------------------------------------------------
# View state: scroll position + current selection
# ---------------------------------------------------------------------
$ws.Range("A11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E19").Select() | Out-Null
